$wb = $excel.ActiveWorkbook

# --- Rename sheets (remove spaces from the 4 tab names) ---
$wsPaises    = $wb.Worksheets.Item(1)   # "Paises de Asia"
$wsTasa      = $wb.Worksheets.Item(2)   # "Grafica Tasa de crecimiento"
$wsHistoria  = $wb.Worksheets.Item(3)   # "Historia de la población de Asi"
$wsGrafHist  = $wb.Worksheets.Item(4)   # "Grafica Historia"

$wsPaises.Name   = "PaisesdeAsia"
$wsTasa.Name     = "GraficaTasadecrecimiento"
$wsHistoria.Name = "Historiadelapoblación de Asi"
$wsGrafHist.Name = "GraficaHistoria"

# --- Fix up the chart series formulas on the "GraficaHistoria" sheet so
#     they reference the renamed "Historiadelapoblación de Asi" sheet ---
$chartObj = $wsGrafHist.ChartObjects().Item(1)
$chart = $chartObj.Chart

$s1 = $chart.SeriesCollection().Item(1)
$s1.Formula = "=SERIES('Historiadelapoblación de Asi'!`$A`$76,'Historiadelapoblación de Asi'!`$A`$2:`$A`$76,'Historiadelapoblación de Asi'!`$A`$2:`$A`$76,1)"

$s2 = $chart.SeriesCollection().Item(2)
$s2.Formula = "=SERIES('Historiadelapoblación de Asi'!`$B`$1,,'Historiadelapoblación de Asi'!`$B`$2:`$B`$76,2)"

# --- Update selection on "PaisesdeAsia" (was D53, now C65) ---
[void]$wsPaises.Range("C65").Select()

# --- Set the page to portrait orientation on "GraficaTasadecrecimiento" ---
$wsTasa.PageSetup.Orientation = 1

# --- Finally activate "GraficaHistoria" and select K19 there; this also
#     moves tabSelected/activeTab to this (4th, index 3) sheet ---
[void]$wsGrafHist.Range("K19").Select()
$wsGrafHist.Activate()
